$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$ws.Activate()

# Remove the "accrualperiodic" row (originally row 33) from the loan product
# input sheet; Excel shifts everything below it up by one row.
[void]$ws.Rows("33:33").Delete()

# Leave the selection on the cell that now occupies row 34 (originally the
# "fundsource" row), matching the post-edit view state.
[void]$ws.Range("A34").Select()
